$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "C"=0.3178392843478832; "D"=0.072975931845086; "E"=0.1152519610375153; "F"=2.663086912965028; "G"=0.002528166556933808; "J"=0.1801976948892259; "L"=0.1004363845456098; "M"=1.795126779606832; "N"=1.770456784104312; "O"=7.650565837941315 }
  3 = @{ "C"=0.3191506642733657; "D"=0.07323142302829844; "E"=0.1165783141418422; "F"=2.64674256577392; "G"=0.00253329629916117; "J"=0.1828342323282115; "L"=0.1002543018401951; "M"=1.669526507629953; "N"=1.641295309442739; "O"=7.583767025111797 }
  4 = @{ "C"=0.3201237289109287; "D"=0.07339847739880945; "E"=0.1174366614359648; "F"=2.638315273754941; "G"=0.002536614178842142; "J"=0.1845423758235372; "L"=0.1001560439589042; "M"=1.592429172857464; "N"=1.562278028786039; "O"=7.547513017343249 }
  5 = @{ "C"=0.320562408484804; "D"=0.07346911763896991; "E"=0.1177974964592072; "F"=2.63528441333608; "G"=0.002538008678222947; "J"=0.1852608790122128; "L"=0.1001194357085247; "M"=1.561019502537576; "N"=1.530154126625064; "O"=7.533932754563978 }
  6 = @{ "C"=0.3206377946399215; "D"=0.07348100238717947; "E"=0.1178580802434017; "F"=2.634805470490974; "G"=0.002538242801178118; "J"=0.1853815390908151; "L"=0.1001135650670442; "M"=1.555804519556489; "N"=1.524824696000024; "O"=7.531749745222044 }
  7 = @{ "C"=0.3201294745271497; "D"=0.07339941968913166; "E"=0.117441483040118; "F"=2.638272766872461; "G"=0.002536632813474708; "J"=0.1845519750737381; "L"=0.1001555363124425; "M"=1.592005534661851; "N"=1.561844480722954; "O"=7.547325041364331 }
  8 = @{ "C"=0.3182565665212564; "D"=0.0730619143608191; "E"=0.1157001618772682; "F"=2.657116941621595; "G"=0.002529900468343256; "J"=0.1810882187743363; "L"=0.1003708081892576; "M"=1.751816839365901; "N"=1.725864252790785; "O"=7.626543305582743 }
  9 = @{ "C"=0.3159189761081507; "D"=0.07248067573396089; "E"=0.1126343243505568; "F"=2.706884802377488; "G"=0.002518026453715672; "J"=0.1750057739647204; "L"=0.100899247335839; "M"=2.065278440517147; "N"=2.049647962062863; "O"=7.819846092821649 }
  10 = @{ "C"=0.3150202932058761; "D"=0.07210253088702601; "E"=0.1105945755317401; "F"=2.751344950762245; "G"=0.002510103189029625; "J"=0.1709716295199426; "L"=0.1013508132440535; "M"=2.295511574237366; "N"=2.288667306981893; "O"=7.985285183608937 }
  11 = @{ "C"=0.3147903178363549; "D"=0.07194106674031264; "E"=0.1097128377549836; "F"=2.773304336917562; "G"=0.00250667059357812; "J"=0.1692311588027184; "L"=0.1015696846022109; "M"=2.40021427016336; "N"=2.397616382568685; "O"=8.065696131948584 }
  12 = @{ "C"=0.3147290339363451; "D"=0.0718814385405393; "E"=0.1093855869730316; "F"=2.781870567931065; "G"=0.002505395306800095; "J"=0.1685857411887426; "L"=0.1016544729611724; "M"=2.439855629538101; "N"=2.438900644764658; "O"=8.096891178801854 }
  13 = @{ "C"=0.314741083554722; "D"=0.07189421321943357; "E"=0.1094557706652022; "F"=2.780014506379416; "G"=0.002505668872244717; "J"=0.168724135114811; "L"=0.1016361279198357; "M"=2.431318520730912; "N"=2.430008171035752; "O"=8.090139556180532 }
  14 = @{ "C"=0.3147847583511094; "D"=0.07193613075200211; "E"=0.1096857814349279; "F"=2.774004054186406; "G"=0.002506565183287404; "J"=0.1691777858591816; "L"=0.1015766221410743; "M"=2.403475750188022; "N"=2.401012331469985; "O"=8.068247607958313 }
  15 = @{ "C"=0.3148148732767169; "D"=0.07196200362600358; "E"=0.1098275351587272; "F"=2.770355164752004; "G"=0.002507117394562; "J"=0.1694574404539555; "L"=0.1015404205958355; "M"=2.386420220221851; "N"=2.383255028525923; "O"=8.0549353297896 }
  16 = @{ "C"=0.3150389286610107; "D"=0.07211329508836428; "E"=0.1106531285268457; "F"=2.749944861823963; "G"=0.002510330961777378; "J"=0.1710872825134437; "L"=0.1013367777886387; "M"=2.288668137211971; "N"=2.281551294532164; "O"=7.980134176854733 }
  17 = @{ "C"=0.3152222452369671; "D"=0.07220880884768377; "E"=0.1111714296316362; "F"=2.737868906191665; "G"=0.002512346271990773; "J"=0.1721114204950558; "L"=0.1012152750111923; "M"=2.228690365887786; "N"=2.219212503123686; "O"=7.935568674099386 }
  18 = @{ "C"=0.3153445121445202; "D"=0.07226473958651614; "E"=0.1114738865275167; "F"=2.731086291700962; "G"=0.002513521598023167; "J"=0.17270938671275; "L"=0.1011466576275311; "M"=2.19418993473181; "N"=2.183377615689096; "O"=7.91042043187349 }
  19 = @{ "C"=0.3153887972598426; "D"=0.07228384754817885; "E"=0.1115770390694466; "F"=2.728817797635131; "G"=0.002513922324871032; "J"=0.1729133765657744; "L"=0.1011236435553009; "M"=2.1825082842212; "N"=2.171248214722539; "O"=7.901988778290217 }
  20 = @{ "C"=0.3152009886192957; "D"=0.07219853841840251; "E"=0.1111158059756516; "F"=2.739137516933837; "G"=0.002512130066069943; "J"=0.1720014768702827; "L"=0.101228078179016; "M"=2.235075412179413; "N"=2.225846454953967; "O"=7.940262556846221 }
  21 = @{ "C"=0.3147712290335107; "D"=0.07192377747174206; "E"=0.1096180413165533; "F"=2.77576265558065; "G"=0.00250630124952097; "J"=0.1690441665193818; "L"=0.1015940489095115; "M"=2.411654066136236; "N"=2.409528392518439; "O"=8.074657545388334 }
  22 = @{ "C"=0.3146407954066177; "D"=0.07175303311254666; "E"=0.1086778963572699; "F"=2.801161122954596; "G"=0.002502634889816499; "J"=0.1671910408666815; "L"=0.1018443317907654; "M"=2.527014739282123; "N"=2.529734406474518; "O"=8.166837892122942 }
  23 = @{ "C"=0.3146966164077583; "D"=0.0718433558717404; "E"=0.1091761231038415; "F"=2.787471289385707; "G"=0.002504578644323269; "J"=0.1681727858924837; "L"=0.1017097446340465; "M"=2.465449475024343; "N"=2.46556484224908; "O"=8.117240449054066 }
  24 = @{ "C"=0.3152105461763597; "D"=0.07220317850654823; "E"=0.1111409394803569; "F"=2.738563479682625; "G"=0.002512227760848063; "J"=0.1720511538049916; "L"=0.1012222860139964; "M"=2.232188789884503; "N"=2.222847230984144; "O"=7.938138979227858 }
  25 = @{ "C"=0.3164079263562058; "D"=0.07262931251183247; "E"=0.1134263611799153; "F"=2.69204081985275; "G"=0.0025210974476092; "J"=0.1765750439417122; "L"=0.100899247335839; "M"=1.980482630772514; "N"=2.049647962062863; "O"=7.819846092821649 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $data[$row][$col]
    }
}